$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 2 (Beta) values F2:N2
$ws.Range("F2").Value = 287.9014347606077
$ws.Range("G2").Value = 13.95553867282598
$ws.Range("H2").Value = 568.0233158455277
$ws.Range("I2").Value = 1.39501333473321
$ws.Range("J2").Value = 0.4961215955190442
$ws.Range("K2").Value = 2.59351501783817
$ws.Range("L2").Value = 0.2175637707877623
$ws.Range("M2").Value = 0.02794169764079556
$ws.Range("N2").Value = 0.4585767591284192

# Update existing row 3 (Gamma) values F3:N3
$ws.Range("F3").Value = 0.02040894916710109
$ws.Range("G3").Value = 0.01260212982147101
$ws.Range("H3").Value = 0.02764727434298727
$ws.Range("I3").Value = 0.01901252379400226
$ws.Range("J3").Value = 0.01165282387365212
$ws.Range("K3").Value = 0.02584029357672211
$ws.Range("L3").Value = 0.02032407638856586
$ws.Range("M3").Value = 0.01253058357237997
$ws.Range("N3").Value = 0.0275564169722681

# Add new row 4 (Beta + Gamma)
$ws.Range("A4").Value = 2
$ws.Range("A2").Copy()
$ws.Range("A4").PasteSpecial(-4122)

$ws.Range("B4").Value = "Beta + Gamma"
$ws.Range("C4").Value = 34.28284193674125
$ws.Range("D4").Value = 3.855929217186839
$ws.Range("E4").Value = 0.2699725673075166
$ws.Range("F4").Value = 287.9218437097748
$ws.Range("G4").Value = 13.96814080264745
$ws.Range("H4").Value = 568.0509631198706
$ws.Range("I4").Value = 1.414025858527212
$ws.Range("J4").Value = 0.5077744193926964
$ws.Range("K4").Value = 2.619355311414892
$ws.Range("L4").Value = 0.2378878471763282
$ws.Range("M4").Value = 0.04047228121317553
$ws.Range("N4").Value = 0.4861331761006873
